$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F, shifting the old "District" column (F) to G
$ws.Columns("F").Insert()

# Set header for new Address column (F2) and existing District header moved to G2
$ws.Range('F2').Value = 'Address'

# Fill in Address values for each row (empty string where address data is absent)
$ws.Range('F3').Value = 'Shree Shankaraling High School Hombal'
$ws.Range('F4').Value = 'J T H S DambalMundaragi'
$ws.Range('F5').Value = 'Sanjay High School HolealurRon'
$ws.Range('F6').Value = 'K H P G G H S Hulkoti'
$ws.Range('F7').Value = 'S M H S Balaganur'
$ws.Range('F8').Value = 'S H S Comp P U CollegePetha – AlurMundaragi'
$ws.Range('F9').Value = 'G H S Shagoti'
$ws.Range('F10').Value = 'G H S (R M S A) Nagavi'
$ws.Range('F11').Value = 'G H S Ron'
$ws.Range('F12').Value = 'M G M High School NyayadaguntePavagada'
$ws.Range('F13').Value = 'G H S B S BeleriRon'
$ws.Range('F14').Value = 'S R A G A H S Hulkoti'
$ws.Range('F16').Value = 'P R H S PalavalliPavagada'
$ws.Range('F17').Value = 'N V A G H S MagadiShirahatti'
$ws.Range('F18').Value = 'G H S (RMSA) HarogeriMundaragi'
$ws.Range('F19').Value = 'M S Dambal Girls High School Mundargi'
$ws.Range('F20').Value = 'G H S Shirahatti'
$ws.Range('F21').Value = 'S R R High School K T HallyPavagada'
$ws.Range('F22').Value = 'Govt. High School Harlapur'
$ws.Range('F23').Value = 'G H S HirekoppaNaragund'
$ws.Range('F24').Value = 'V F Patil High School Ron'
$ws.Range('F25').Value = 'G H S Siddaling Nagar'
$ws.Range('F26').Value = 'S V V P Girls High School Pavagada'
$ws.Range('F27').Value = 'Nethra Vidya Peeta High School GujjanaduPavagada'
$ws.Range('F29').Value = 'H C E S High School Chincholi'
$ws.Range('F30').Value = 'G H S HadaliNaragund'
$ws.Range('F31').Value = 'G H S KoganurShirahatti'
$ws.Range('F32').Value = 'G H S Janili – ShirurMundaragi'
$ws.Range('F33').Value = 'Akarashni High School NeelammanahallyPavagada'
$ws.Range('F34').Value = 'S K B H S MadalageriRon'
$ws.Range('F35').Value = 'Sri K V Shantagirimath High SchoolRon'
$ws.Range('F36').Value = 'M H S ShanthgiriRon'
$ws.Range('F37').Value = 'G H S BaradurMundaragi'
$ws.Range('F38').Value = 'G H S MenasagiRon'
$ws.Range('F39').Value = 'J A H S Mundargi'
$ws.Range('F40').Value = 'Smt. PSBD G H S LakshmishwarShirahatti'
$ws.Range('F42').Value = 'G H S RamagiriShirahatti'
$ws.Range('F43').Value = 'G G H S Lakkundi'
$ws.Range('F44').Value = 'S J F H S BellattiShirahatti'
$ws.Range('F45').Value = 'G H S ChikkamagundNaragund'
$ws.Range('F46').Value = 'S S High School LakshmeshwarShirahatti'
$ws.Range('F47').Value = 'S A D S G H S JakkaliRon'
$ws.Range('F49').Value = 'S T R High SchoolRangasamudraPavagada'
$ws.Range('F50').Value = 'G H S KurahattiRon'

# Touch column H so the sheet's used-range/dimension extends to H50,
# matching the target dimension ref "A1:H50" (no visible content is added;
# the pattern is set to the already-default "none" so no new style/content
# is actually introduced).
$ws.Range('H50').Interior.Pattern = -4142
